$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Title heading + bold "Play Aurora..." near end share identical old text,
# so Find/Replace with Wrap=1 (wdFindContinue) will hit all occurrences.
Replace-Text "Play Aurora Beast Hunter Free Today" "Play Aurora Beast Hunter Free - Exciting Wilds, Bonus Round, and Free Spins"

# "What we like" bullet list
# NOTE: order matters here. The new text for the "High variance..." bullet
# contains "Ability to choose from 3 different free spin features" as a
# substring, so the old "Ability to choose..." bullet must be replaced
# FIRST, before that substring is (re)introduced elsewhere in the document.
Replace-Text "Expanding and stacked wilds increase winning potential" "Exciting gameplay features with wilds, bonus round, and free spins"
Replace-Text "Ability to choose from 3 different free spin features" "High variance game with a maximum win of 5000x your bet"
Replace-Text "High variance and 96.27% RTP offer big wins" "Ability to choose from 3 different free spin features with varying variance"
Replace-Text "Excellent visual design and animations" "Excellent retro animated design and bold graphics"

# "What we don't like" bullet list
Replace-Text "Success rate may discourage casual players" "Success rate of about 1 in every 5 spins may not appeal to some players"
Replace-Text "Bonus rounds may be difficult to trigger" "Limited bet range with a maximum bet of €50"

# Closing italic summary paragraph
Replace-Text "Experience big wins and exciting features in Aurora Beast Hunter, a 5-reel, 40-payline slot game from Just for the Win and Microgaming. Play free now." "Read our review of Aurora Beast Hunter, an exciting slot game with wilds, bonus round, and free spins. Play for free now."
